$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 127; existing rows 127..192 shift down to 128..193
$ws.Rows.Item(127).Insert()

# Populate the newly inserted row 127 with the new record
$ws.Range("A127").Value = 1
$ws.Range("B127").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C127").Value = "Arica y Parinacota"
$ws.Range("D127").Value = 44596
$ws.Range("E127").Value = 15
$ws.Range("F127").Value = "Fruta"
$ws.Range("G127").Value = 100102
$ws.Range("H127").Value = "Cítricos"
$ws.Range("I127").Value = 100102003
$ws.Range("J127").Value = "Limón"
$ws.Range("K127").Value = "Tahití"
$ws.Range("L127").Value = "Primera"
$ws.Range("M127").Value = 300
$ws.Range("N127").Value = 35000
$ws.Range("O127").Value = 36000
$ws.Range("P127").Value = 35500
$ws.Range("Q127").Value = '$/caja 24 kilos'
$ws.Range("R127").Value = "Perú"
$ws.Range("S127").Value = 1479
$ws.Range("T127").Value = 24

# Match the date format used by the rest of column D
$ws.Range("D127").NumberFormat = $ws.Range("D128").NumberFormat
